# Update pin assignment: swap "nord_"/"south_" labels between the two
# connector groups (rows 2-5 and rows 7-10) in column A (and the
# complementary G column labels), then move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "(nord_3V3)"
$ws.Range("A3").Value = "nord_channel"
$ws.Range("A4").Value = "nord_up"
$ws.Range("G4").Value = "nord_GND"
$ws.Range("A5").Value = "nord_down"
$ws.Range("G5").Value = "nord_my"

$ws.Range("A7").Value = "south_up"
$ws.Range("A8").Value = "south_my"
$ws.Range("G8").Value = "south_GND"
$ws.Range("A9").Value = "south_down"
$ws.Range("G9").Value = "south_channel"
$ws.Range("A10").Value = "(south_3V3)"

$ws.Range("A8").Select()
